$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cell values
$ws.Range("B9").Value = "Postman2UPDATE2UPDATE2"
$ws.Range("B12").Value = "Postman5PT3"

# Append new rows 18-21, using row 17's formatting as the template
$ws.Range("A17:I17").Copy($ws.Range("A18:I18"))
$ws.Range("A17:I17").Copy($ws.Range("A19:I19"))
$ws.Range("A17:I17").Copy($ws.Range("A20:I20"))
$ws.Range("A17:I17").Copy($ws.Range("A21:I21"))

$newRows = @(
    @{ Row = 18; A = 16; B = "Postman10";          C = 44958.59305555555; D = "AGAIN I sent yet ANOTHER email via postman"; E = "marleevaughn@outlook.com"; F = "Marlee Vaughn"; G = "duanevaughn@hotmail.com";  H = "Duane Vaughn";       I = $true },
    @{ Row = 19; A = 17; B = "Postman10";          C = 44958.59305555555; D = "AGAIN I sent yet ANOTHER email via postman"; E = "marleevaughn@outlook.com"; F = "Marlee Vaughn"; G = "duanevaughn@hotmail.com";  H = "Duane Vaughn";       I = $true },
    @{ Row = 20; A = 18; B = "Postman11";          C = 44958.59305555555; D = "AGAIN I sent yet ANOTHER email via postman"; E = "marleevaughn@outlook.com"; F = "Marlee Vaughn"; G = "duanevaughn@hotmail.com";  H = "Duane Vaughn";       I = $true },
    @{ Row = 21; A = 19; B = "Saving new email 3"; C = 44958.59305555555; D = "Hope this works AGAIN!";                     E = "duanevaughn@hotmail.com";  F = "Duane Vaughn";  G = "{No Recipient Email}";    H = "{No Recipient Name}"; I = $true }
)

foreach ($rowData in $newRows) {
    $r = $rowData.Row
    $ws.Cells.Item($r, 1).Value = $rowData.A
    $ws.Cells.Item($r, 2).Value = $rowData.B
    $ws.Cells.Item($r, 3).Value = $rowData.C
    $ws.Cells.Item($r, 4).Value = $rowData.D
    $ws.Cells.Item($r, 5).Value = $rowData.E
    $ws.Cells.Item($r, 6).Value = $rowData.F
    $ws.Cells.Item($r, 7).Value = $rowData.G
    $ws.Cells.Item($r, 8).Value = $rowData.H
    $ws.Cells.Item($r, 9).Value = $rowData.I
}
